$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 6221071.428571429
$ws.Range("C3").Value = 3915000
$ws.Range("C4").Value = 3163928.571428571
$ws.Range("C5").Value = 6630000
$ws.Range("C6").Value = 3185000
$ws.Range("C7").Value = 23115000
